$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 201, shifting existing rows 201:246 down to 202:247
$ws.Rows(201).Insert()

# Populate the new row 201 with the new data point
$ws.Range("A201").Value = 5
$ws.Range("B201").Value = "Macroferia Regional de Talca"
$ws.Range("C201").Value = "Maule"
$ws.Range("D201").Value = 44551
$ws.Range("E201").Value = 7
$ws.Range("F201").Value = 100112023
$ws.Range("G201").Value = "Brócoli"
$ws.Range("H201").Value = "Sin especificar"
$ws.Range("I201").Value = "Primera"
$ws.Range("J201").Value = 3000
$ws.Range("K201").Value = 700
$ws.Range("L201").Value = 700
$ws.Range("M201").Value = 700
$ws.Range("N201").Value = "$/unidad"
$ws.Range("O201").Value = "Región del Maule"
$ws.Range("P201").Value = 700
$ws.Range("Q201").Value = 1
$ws.Range("R201").Value = "Hortaliza"
